# Updates cryptos list values (Price & Volume(1h) columns, plus a few
# coin name/link swaps) to match the latest scrape, per the commit:
# "Updated cryptos list on Tue Feb 27 03:58:49 UTC 2024 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (B Coin, C Link, D Price, E Volume(1h))
# Only rows whose values actually changed are listed; all rows 2-51 end up
# consistent with the final state described by the diff.
$rows = @{
    2  = @('Bitcoin',                      'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc',                       '55.913.76',  '  +8.54%  ')
    3  = @('Ethereum',                     'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth',                      '3.224.61',   '  +3.99%  ')
    4  = @('TetherUSD',                    'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt',                    '1.00',       '  +0.18%  ')
    5  = @('BNB',                          'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb',                           '400.72',     '  +4.29%  ')
    6  = @('Solana',                       'https://coinranking.com/coin/zNZHO_Sjf+solana-sol',                            '110.77',     '  +7.32%  ')
    7  = @('XRP',                          'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp',                           '0.555',      '  +2.75%  ')
    8  = @('USDC',                         'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc',                         '0.999',      '  -0.05%  ')
    9  = @('Cardano',                      'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada',                       '0.623',      '  +6.62%  ')
    10 = @('Avalanche',                    'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax',                        '39.30',      '  +6.39%  ')
    11 = @('Dogecoin',                     'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge',                     '0.0897',     '  +4.81%  ')
    12 = @('TRON',                         'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx',                          '0.140',      '  +2.00%  ')
    13 = @('WrappedliquidstakedEther2.0',  'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth',     '3.751.18',   '  +4.58%  ')
    14 = @('Polkadot',                     'https://coinranking.com/coin/25W7FG7om+polkadot-dot',                          '8.07',       '  +2.89%  ')
    15 = @('Chainlink',                    'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link',                    '19.05',      '  +2.20%  ')
    16 = @('WrappedEther',                 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth',                 '3.236.94',   '  +4.43%  ')
    17 = @('Polygon',                      'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic',                      '1.06',       '  +6.72%  ')
    18 = @('Uniswap',                      'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni',                          '10.61',      '  -4.37%  ')
    19 = @('WrappedBTC',                   'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc',                    '55.969.72',  '  +8.56%  ')
    20 = @('ImmutableX',                   'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx',                        '3.39',       '  +1.60%  ')
    21 = @('ShibaInu',                     'https://coinranking.com/coin/xz24e0BjL+shibainu-shib',                         '0.0000102',  '  +5.91%  ')
    22 = @('InternetComputer(DFINITY)',    'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp',           '13.04',      '  +5.18%  ')
    23 = @('BitcoinCash',                  'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch',                   '301.00',     '  +13.06%  ')
    24 = @('Litecoin',                     'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc',                      '74.93',      '  +7.04%  ')
    25 = @('PancakeSwap',                  'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake',                      '3.24',       '  +2.96%  ')
    26 = @('Filecoin',                     'https://coinranking.com/coin/ymQub4fuB+filecoin-fil',                          '8.17',       '  +0.30%  ')
    27 = @('EthereumClassic',              'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc',               '28.19',      '  +4.11%  ')
    28 = @('RenderToken',                  'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr',                  '7.42',       '  +2.21%  ')
    29 = @('Kaspa',                        'https://coinranking.com/coin/V8GxkwWow+kaspa-kas',                             '0.172',      '  +2.69%  ')
    30 = @('Dai',                          'https://coinranking.com/coin/MoTuySvg7+dai-dai',                               '0.998',      '  -0.26%  ')
    31 = @('Hedera',                       'https://coinranking.com/coin/jad286TjB+hedera-hbar',                           '0.112',      '  +4.34%  ')
    32 = @('Cosmos',                       'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom',                       '11.35',      '  +9.72%  ')
    33 = @('VeChain',                      'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet',                       '0.0491',     '  +5.31%  ')
    34 = @('InjectiveProtocol',            'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj',                 '36.60',      '  +3.31%  ')
    35 = @('Toncoin',                      'https://coinranking.com/coin/67YlI0K1b+toncoin-ton',                           '2.11',       '  +2.03%  ')
    36 = @('OKB',                          'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb',                           '51.39',      '  +2.18%  ')
    37 = @('LidoDAOToken',                 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo',                      '3.55',       '  +5.66%  ')
    38 = @('FirstDigitalUSD',              'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd',                 '1.00',       '  +0.21%  ')
    39 = @('Stacks',                       'https://coinranking.com/coin/mMPrMcB7+stacks-stx',                             '3.08',       '  +22.36%  ')
    40 = @('ARBITRUM',                     'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb',                          '1.92',       '  +1.41%  ')
    41 = @('NEARProtocol',                 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near',                     '4.02',       '  +9.83%  ')
    42 = @('Monero',                       'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr',                        '131.79',     '  +2.22%  ')
    43 = @('Celestia',                     'https://coinranking.com/coin/YQcD0lBl7+celestia-tia',                         '17.10',      '  +3.15%  ')
    44 = @('Stellar',                      'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm',                       '0.119',      '  +2.71%  ')
    45 = @('TheGraph',                     'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt',                          '0.283',      '  -5.94%  ')
    46 = @('EnergySwap',                   'https://coinranking.com/coin/SbWqqTui-+energyswap-ens',                        '22.26',      '  -0.66%  ')
    47 = @('ThetaToken',                   'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta',                  '2.11',       '  +42.49%  ')
    48 = @('Maker',                        'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr',                         '2.156.68',   '  +4.67%  ')
    49 = @('WEMIXToken',                   'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix',                      '2.09',       '  +0.91%  ')
    50 = @('ApeXProtocol',                 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex',                     '2.46',       '  -0.03%  ')
    51 = @('BEAM',                         'https://coinranking.com/coin/cYYMfXF4u+beam-beam',                             '0.0357',     '  +7.52%  ')
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
}
